$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A54").Value = "Francesco Passuello"
$ws.Range("B54").Value = "Gabriel Melis | demobusters"
$ws.Range("C54").Value = "Mattia Baldessarini | SHARK ATTACK"
$ws.Range("D54").Value = "Federico Mortillaro | Clitoriders"
$ws.Range("E54").Value = "Mattia Tezzele | U.SGUARNA"
$ws.Range("F54").Value = "Thomas Pontillo | Gli Introvabili"
